$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.022.21"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.67%  '

# Row 3
$ws.Range("D3").Value = "'2.408.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.92%  '

# Row 4
$ws.Range("E4").Value = '  +0.33%  '

# Row 5
$ws.Range("D5").Value = "'560.25"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.30%  '

# Row 6
$ws.Range("D6").Value = "'164.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.68%  '

# Row 7
$ws.Range("E7").Value = '  +0.11%  '

# Row 8
$ws.Range("D8").Value = "'0.509"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.24%  '

# Row 9
$ws.Range("D9").Value = "'0.167"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.81%  '

# Row 10
$ws.Range("D10").Value = "'2.405.77"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.98%  '

# Row 11
$ws.Range("E11").Value = '  -1.93%  '

# Row 12
$ws.Range("D12").Value = "'0.330"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.52%  '

# Row 13
$ws.Range("D13").Value = "'4.65"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.68%  '

# Row 14
$ws.Range("B14").Value = 'WrappedBTC'
$ws.Range("C14").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D14").Value = "'69.196.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.04%  '

# Row 15
$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").Value = "'0.0000177"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.99%  '

# Row 16
$ws.Range("D16").Value = "'2.873.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.10%  '

# Row 17
$ws.Range("D17").Value = "'23.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.52%  '

# Row 18
$ws.Range("D18").Value = "'2.439.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.63%  '

# Row 19
$ws.Range("D19").Value = "'10.67"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.84%  '

# Row 20
$ws.Range("D20").Value = "'337.29"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.19%  '

# Row 21
$ws.Range("D21").Value = "'7.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.66%  '

# Row 22
$ws.Range("D22").Value = "'3.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.38%  '

# Row 23
$ws.Range("D23").Value = "'1.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.78%  '

# Row 24
$ws.Range("D24").Value = "'1.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.05%  '

# Row 25
$ws.Range("D25").Value = "'65.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.03%  '

# Row 26
$ws.Range("D26").Value = "'3.76"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.82%  '

# Row 27
$ws.Range("D27").Value = "'2.568.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.46%  '

# Row 28
$ws.Range("D28").Value = "'8.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.98%  '

# Row 29
$ws.Range("D29").Value = "'1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.20%  '

# Row 30
$ws.Range("D30").Value = "'0.0₃0840"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.47%  '

# Row 31
$ws.Range("D31").Value = "'7.26"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.89%  '

# Row 32
$ws.Range("D32").Value = "'1.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.69%  '

# Row 33
$ws.Range("D33").Value = "'449.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.50%  '

# Row 34
$ws.Range("E34").Value = '  +0.29%  '

# Row 35
$ws.Range("D35").Value = "'1.60"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.36%  '

# Row 36
$ws.Range("D36").Value = "'159.31"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.61%  '

# Row 37
$ws.Range("E37").Value = '  +0.74%  '

# Row 38
$ws.Range("B38").Value = 'USDe'
$ws.Range("C38").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.04%  '

# Row 39
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = "'0.110"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.30%  '

# Row 40
$ws.Range("D40").Value = "'18.02"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.09%  '

# Row 41
$ws.Range("D41").Value = "'0.301"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.31%  '

# Row 42
$ws.Range("D42").Value = "'37.87"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.33%  '

# Row 43
$ws.Range("D43").Value = "'1.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.69%  '

# Row 44
$ws.Range("D44").Value = "'4.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.83%  '

# Row 45
$ws.Range("D45").Value = "'1.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.12%  '

# Row 46
$ws.Range("D46").Value = "'2.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.10%  '

# Row 47
$ws.Range("D47").Value = "'133.22"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.87%  '

# Row 48
$ws.Range("D48").Value = "'3.34"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.63%  '

# Row 49
$ws.Range("D49").Value = "'0.0723"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.64%  '

# Row 50
$ws.Range("D50").Value = "'0.484"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.27%  '

# Row 51
$ws.Range("D51").Value = "'0.556"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.12%  '
